$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (A97) duplicating the value already present in A96
# (using Value() to properly read the getter value from this COM runtime)
$lastValue = $ws.Range("A96").Value()
$ws.Range("A97").Value = $lastValue
